# Applies the "double underscore as field separator for flattened fields" change
# plus related legend/label wording tweaks and column width adjustments.

$wb = $excel.ActiveWorkbook

$wsSchema = $wb.Worksheets.Item("semantic_aspect_model_schema")
$wsDescription = $wb.Worksheets.Item("description")

# ---------------------------------------------------------------------------
# 1. semantic_aspect_model_schema: adjust column widths for columns 8-20
#
# Note: the Excel column-width engine here stores widths quantized to a
# Calibri-11 (MDW=7px) pixel grid, so the requested target widths are
# approached via the input value that the engine snaps closest to the
# intended target (the same quantization genuine Excel applies to
# ColumnWidth assignments).
# ---------------------------------------------------------------------------
$colWidths = @{
    8  = 48.333333333333336   # -> stored width 49.166666... (target 49.2)
    9  = 42.333333333333336   # -> stored width 43.166666... (target 43.2)
    10 = 31.5                 # -> stored width 32.333333... (target 32.4)
    11 = 43.5                 # -> stored width 44.333333... (target 44.4)
    12 = 46.0                 # -> stored width 46.833333... (target 46.8)
    13 = 37.5                 # -> stored width 38.333333... (target 38.4)
    14 = 49.166666666666664   # -> stored width 50            (target 50)
    15 = 40.0                 # -> stored width 40.833333... (target 40.8)
    16 = 41.166666666666664   # -> stored width 42            (target 42)
    17 = 35.166666666666664   # -> stored width 36            (target 36)
    18 = 40.0                 # -> stored width 40.833333... (target 40.8)
    19 = 40.0                 # -> stored width 40.833333... (target 40.8)
    20 = 38.833333333333336   # -> stored width 39.666666... (target 39.6)
}

foreach ($col in $colWidths.Keys) {
    $wsSchema.Columns.Item($col).ColumnWidth = $colWidths[$col]
}

# ---------------------------------------------------------------------------
# 2. semantic_aspect_model_schema: row 1 header renames (single "_" -> "__")
# ---------------------------------------------------------------------------
$headerRenames = @{
    "D1" = "cxPCFValues__sharesPerStage__stage1RawMaterialAquisitionPreProcessing__stageValue"
    "E1" = "cxPCFValues__sharesPerStage__stage2MainProductProduction__stageValue"
    "F1" = "cxPCFValues__sharesPerStage__stage3Distribution__stageValue"
    "G1" = "cxPCFValues__sharesPerStage__stage4EndOfLifeRecycling__stageValue"
    "H1" = "cxPCFValues__sharesPerStage__declaredUnit"
    "I1" = "cxPCFValues__rawData[0]__description"
    "J1" = "cxPCFValues__rawData[0]__id"
    "K1" = "cxPCFValues__rawData[0]__data[0]__key"
    "L1" = "cxPCFValues__rawData[0]__data[0]__value"
    "M1" = "cxPCFValues__pcfPerformanceClass"
    "N1" = "cxPCFValues__webLinkToCarbonFootPrintStudy"
    "O1" = "dataQualityRating__coveragePercent"
    "P1" = "dataQualityRating__technologicalDQR"
    "Q1" = "dataQualityRating__temporalDQR"
    "R1" = "dataQualityRating__geographicalDQR"
    "S1" = "dataQualityRating__completenessDQR"
    "T1" = "dataQualityRating__reliabilityDQR"
    "U1" = "secondaryEmissionFactorSources[0]__secondaryEmissionFactorSource"
}

foreach ($cellRef in $headerRenames.Keys) {
    $wsSchema.Range($cellRef).Value = $headerRenames[$cellRef]
}

# ---------------------------------------------------------------------------
# 3. description sheet: legend numbering fix
# ---------------------------------------------------------------------------
$wsDescription.Range("A3").Value = "1. Columns highlighted in olive green are digital twin fields."

# ---------------------------------------------------------------------------
# 4. description sheet: "Digital Twin Field:" -> "Digital Twin Field Name:"
# ---------------------------------------------------------------------------
$wsDescription.Range("B5").Value = "Digital Twin Field Name: id"
$wsDescription.Range("B6").Value = "Digital Twin Field Name: manufacturerPartId"
$wsDescription.Range("B7").Value = "Digital Twin Field Name: partInstanceId"

# ---------------------------------------------------------------------------
# 5. description sheet: column A field names (single "_" -> "__")
# ---------------------------------------------------------------------------
$descriptionRenames = @{
    "A8"  = "cxPCFValues__sharesPerStage__stage1RawMaterialAquisitionPreProcessing__stageValue"
    "A9"  = "cxPCFValues__sharesPerStage__stage2MainProductProduction__stageValue"
    "A10" = "cxPCFValues__sharesPerStage__stage3Distribution__stageValue"
    "A11" = "cxPCFValues__sharesPerStage__stage4EndOfLifeRecycling__stageValue"
    "A12" = "cxPCFValues__sharesPerStage__declaredUnit"
    "A13" = "cxPCFValues__rawData[0]__description"
    "A14" = "cxPCFValues__rawData[0]__id"
    "A15" = "cxPCFValues__rawData[0]__data[0]__key"
    "A16" = "cxPCFValues__rawData[0]__data[0]__value"
    "A17" = "cxPCFValues__pcfPerformanceClass"
    "A18" = "cxPCFValues__webLinkToCarbonFootPrintStudy"
    "A19" = "dataQualityRating__coveragePercent"
    "A20" = "dataQualityRating__technologicalDQR"
    "A21" = "dataQualityRating__temporalDQR"
    "A22" = "dataQualityRating__geographicalDQR"
    "A23" = "dataQualityRating__completenessDQR"
    "A24" = "dataQualityRating__reliabilityDQR"
    "A25" = "secondaryEmissionFactorSources[0]__secondaryEmissionFactorSource"
}

foreach ($cellRef in $descriptionRenames.Keys) {
    $wsDescription.Range($cellRef).Value = $descriptionRenames[$cellRef]
}

$wb.Save()
